# Auto-generated Excel COM-interop script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.409.78"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.526.82"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'315.04"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").Value = "'94.12"
$ws.Range("E6").Value = "  -5.22%  "
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").Value = "'35.55"
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("D12").Value = "'7.60"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "2.913.30"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.524.94"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.42"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "'0.844"
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").Value = "42.477.73"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "'12.85"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").Value = "'6.56"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").Value = "'70.57"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "'249.51"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "'2.01"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("D26").Value = "'26.62"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").Value = "'39.01"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'10.11"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "'5.92"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("D32").Value = "'155.83"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "'19.26"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'2.12"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.30"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'0.0781"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -4.50%  "
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").Value = "'23.71"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").Value = "'2.34"
$ws.Range("E41").Value = "  +10.35%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'3.80"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'3.31"
$ws.Range("E44").Value = "  -4.56%  "
$ws.Range("D45").Value = "'0.0299"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "2.019.54"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "'84.24"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").Value = "'8.81"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").Value = "2.765.40"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'102.03"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "'72.59"
$ws.Range("E51").Value = "  -1.62%  "
